# docs/artifacts/rules/dms-addition-svein-harald.xlsx
# Insert a new "Immutable" column (H) into the Properties and LastProperties
# sheets, right after the "Nullable" column, pushing the remaining columns
# (Is List..Property (linage)) one to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Properties sheet
# ---------------------------------------------------------------------
$wsProps = $wb.Worksheets.Item("Properties")

# Insert a blank column before the old "Is List" column (H); this shifts
# H:P -> I:Q, extends the merged title range and the column-width groups
# automatically.
$wsProps.Columns("H").Insert()

# New column header
$wsProps.Range("H2").Value = "Immutable"

# Rows that represent real properties (i.e. the "Nullable" column in G is
# populated) get an explicit FALSE default for the new "Immutable" column.
$propsFalseRows = @(3,5,7,8,9,10,11,12)
foreach ($r in $propsFalseRows) {
    $wsProps.Range("H$r").Value = $false
}

# ---------------------------------------------------------------------
# LastProperties sheet
# ---------------------------------------------------------------------
$wsLastProps = $wb.Worksheets.Item("LastProperties")

$wsLastProps.Columns("H").Insert()

$wsLastProps.Range("H2").Value = "Immutable"

$lastPropsFalseRows = @(3,5,7,9,10,11,12,14,15,16,18,19,20,21,23,26,27,31,32,36,37,38,39,44,45,46,47,48,49,50,52,54,55,58,59,62,63,64,65,66)
foreach ($r in $lastPropsFalseRows) {
    $wsLastProps.Range("H$r").Value = $false
}

# ---------------------------------------------------------------------
# View-state: user ends up with "Properties" as the active/selected tab
# with G24 selected, and "LastProperties" with H67 selected (no longer the
# active tab).
# ---------------------------------------------------------------------
$wsLastProps.Range("H67").Select()
$wsProps.Activate()
$wsProps.Range("G24").Select()
